# "se agregó el reseteo de Sec cuando cambia nro de bobina"
# (added resetting of Sec when the coil number changes)
#
# Append new coil-reading rows (177-189) under the existing data table.
# Every column in this sheet is stored as text (Ancho, Diametro, Gramaje,
# Peso, Bobina Nro, Orden de Fabricacion, CodCal all look numeric but
# carry leading zeros / are treated as literal codes) except "Sec" (F),
# which is the value being exercised by this change: it is a real number
# that starts back at 1 whenever "Bobina Nro" (E) changes, and only
# increments (here, to 2) when the same Bobina Nro repeats in the next
# reading (rows 180/181).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A="12";  B="23";  C="23";  D="23";  E="23";   F=1;   G="32";    H="2025-03-06 15:43"; I="B"; J="03"; K="L.BLANCO" },
    @{ A="455"; B="120"; C="130"; D="354"; E="4556";  F=1;   G="54665"; H="2025-03-06 15:44"; I="B"; J="02"; K="COVERING" },
    @{ A="144"; B="120"; C="130"; D="545"; E="4545";  F=1;   G="5444";  H="2025-03-06 15:46"; I="D"; J="02"; K="COVERING" },
    @{ A="344"; B="455"; C="545"; D="455"; E="4532";  F=1;   G="23444"; H="2025-03-06 15:47"; I="B"; J="02"; K="COVERING" },
    @{ A="344"; B="455"; C="545"; D="455"; E="4532";  F="2"; G="23444"; H="2025-03-06 15:47"; I="B"; J="02"; K="COVERING" },
    @{ A="122"; B="120"; C="130"; D="454"; E="458";   F=1;   G="8785";  H="2025-03-06 15:51"; I="C"; J="03"; K="L.BLANCO" },
    @{ A="150"; B="120"; C="130"; D="454"; E="4548";  F="1"; G="8777";  H="2025-03-06 15:53"; I="C"; J="04"; K="CART.GRIS" },
    @{ A="120"; B="120"; C="150"; D="554"; E="4545";  F=1;   G="5455";  H="2025-03-06 15:54"; I="B"; J="06"; K="LINER PER" },
    @{ A="160"; B="120"; C="150"; D="545"; E="5448";  F=1;   G="8885";  H="2025-03-06 15:55"; I="A"; J="04"; K="CART.GRIS" },
    @{ A="150"; B="120"; C="150"; D="545"; E="4555";  F="1"; G="545";   H="2025-03-06 15:56"; I="C"; J="02"; K="COVERING" },
    @{ A="160"; B="120"; C="150"; D="540"; E="4488";  F=1;   G="88652"; H="2025-03-06 15:59"; I="B"; J="02"; K="COVERING" },
    @{ A="160"; B="120"; C="160"; D="548"; E="4589";  F=1;   G="98788"; H="2025-03-06 16:01"; I="B"; J="03"; K="L.BLANCO" },
    @{ A="150"; B="120"; C="130"; D="545"; E="5488";  F=1;   G="87888"; H="2025-03-06 16:02"; I="B"; J="03"; K="L.BLANCO" }
)

$startRow = 177
$lastRow = $startRow + $rows.Count - 1

# Pre-format the numeric-looking columns as Text so digit strings (and
# codes with leading zeros, like CodCal) are kept verbatim instead of
# being coerced into numbers.
$ws.Range("A$($startRow):E$($lastRow)").NumberFormat = "@"
$ws.Range("G$($startRow):G$($lastRow)").NumberFormat = "@"
$ws.Range("J$($startRow):J$($lastRow)").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    if ($row.F -is [string]) {
        # These particular "Sec" readings were captured as text.
        $ws.Cells.Item($r, 6).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 6).Value = $row.F

    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
}
